# The "Requisitos:" section (row 23, column A) is followed by two entries,
# one per requirement, in columns B and C of rows 24 and 25:
#   Row 24: LOB1053 -  Física III  (Requisito)
#   Row 25: LOM3254 -  Laboratório de Circuitos Elétricos  (Indicação de Conjunto)
#
# The change reorders these two entries so the LOM3254 (Indicação de Conjunto)
# entry comes first, followed by the LOB1053 (Requisito) entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$labText = "LOM3254 -  Laboratório de Circuitos Elétricos  (Indicação de Conjunto)`n"
$reqText = "LOB1053 -  Física III  (Requisito)`n"

$ws.Range("B24").Value = $labText
$ws.Range("C24").Value = $labText
$ws.Range("B25").Value = $reqText
$ws.Range("C25").Value = $reqText
